$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit rotates the data in rows 2, 4 and 5 (row 3 stays untouched):
#   new Row2 <- old Row5
#   new Row4 <- old Row2
#   new Row5 <- old Row4
# Columns D, J, K, L, M, P carry the changing values (dates/volumes/prices);
# the rest of the columns happen to hold identical values across these rows.

$cols = @("D", "J", "K", "L", "M", "P")

# Capture original values before overwriting anything.
$orig2 = @{}
$orig4 = @{}
$orig5 = @{}
foreach ($col in $cols) {
    $orig2[$col] = $ws.Range($col + "2").Value2
    $orig4[$col] = $ws.Range($col + "4").Value2
    $orig5[$col] = $ws.Range($col + "5").Value2
}

foreach ($col in $cols) {
    $ws.Range($col + "2").Value2 = $orig5[$col]
    $ws.Range($col + "4").Value2 = $orig2[$col]
    $ws.Range($col + "5").Value2 = $orig4[$col]
}
